$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (capitalize Summer/Winter); column positions stay the same
$ws.Range("B1").Value = "total.crc.Summer"
$ws.Range("C1").Value = "reported.crc.Summer"
$ws.Range("D1").Value = "total.crc.Winter"
$ws.Range("E1").Value = "reported.crc.Winter"

# Add the new 2024 summer data row
$ws.Range("A19").Value = 2024
$ws.Range("B19").Value = 187214
$ws.Range("C19").Value = 81091

# Widen columns B:E to fit new longer headers
$ws.Columns.Item(2).ColumnWidth = 19.053385416666668
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 19.830729166666668
$ws.Columns.Item(5).ColumnWidth = 22.498697916666668

# Scroll the view down and set the new active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$null = $ws.Range("D21").Select()
